$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "current": the "delete history" command should have wiped all
# stored rows, leaving only the header row behind.
# ---------------------------------------------------------------------
$wsCurrent = $wb.Worksheets.Item("current")
$wsCurrent.Rows("2:5").Delete()

# ---------------------------------------------------------------------
# Sheet "forecast": gains a simple placeholder row of data.
# ---------------------------------------------------------------------
$wsForecast = $wb.Worksheets.Item("forecast")
$wsForecast.Range("A1").Value = 1
$wsForecast.Range("B1").Value = 1
$wsForecast.Range("C1").Value = 1
$wsForecast.Range("D1").Value = 1

# ---------------------------------------------------------------------
# Sheet "football": the old row 1 is cleared out, and the history that
# used to live in rows 1-2 now lives in rows 2-5 (row 2 holds the
# "football" summary row, rows 3-5 repeat the latest entry).
# ---------------------------------------------------------------------
$wsFootball = $wb.Worksheets.Item("football")

# Grab the still-valid label strings before we start overwriting cells.
$listMon = $wsFootball.Range("B1").Value()
$listTue = $wsFootball.Range("B2").Value()
$footballLabel = $wsFootball.Range("D1").Value()

$wsFootball.Rows("1:1").ClearContents()

$footballJsonBulgaria = '{"football":[{"stadium":"Fc Rapid Bucuresti","country":"Romania","region":"","tournament":"UEFA Nations League","start":"2022-06-14 19:45","match":"Romania vs Montenegro"},{"stadium":"Inter Turku , Fc Tps","country":"Finland","region":"","tournament":"Finnish Veikkausliga","start":"2022-06-22 16:00","match":"FC Inter Turku vs FC Lahti"},{"stadium":"Drogheda United","country":"Ireland","region":"","tournament":"League of Ireland Premier Division","start":"2022-06-24 19:45","match":"Drogheda vs Sligo Rovers"},{"stadium":"Shelbourne Fc","country":"Ireland","region":"","tournament":"League of Ireland Premier Division","start":"2022-06-24 19:45","match":"Shelbourne vs Dundalk"},{"stadium":"Ucd","country":"Ireland","region":"","tournament":"League of Ireland Premier Division","start":"2022-06-24 19:45","match":"U.C.D vs Derry City"},{"stadium":"Cork City Fc","country":"Ireland","region":"","tournament":"League of Ireland First Division","start":"2022-06-24 19:45","match":"Cork City vs Cobh Ramblers"},{"stadium":"Waterford United Fc","country":"Ireland","region":"","tournament":"League of Ireland First Division","start":"2022-06-24 19:45","match":"Waterford United vs Bray Wanderers"},{"stadium":"Finn Harps","country":"Bulgaria","region":"","tournament":"League of Ireland Premier Division","start":"2022-06-24 20:00","match":"Finn Harps vs St Patricks Athletic"},{"stadium":"Shamrock Rovers(A)","country":"Ireland","region":"","tournament":"League of Ireland Premier Division","start":"2022-06-24 20:00","match":"Shamrock Rovers vs Bohemians"},{"stadium":"Cincinnati Kings","country":"United States of America","region":"","tournament":"American MLS League","start":"2022-06-25 00:30","match":"FC Cincinnati vs Orlando City SC"}],"cricket":[{"stadium":"Sir Vivian Richards Stadium, North Sound, Antigua","country":"West Indies","region":"","tournament":"West Indies vs Bangladesh Test Series 2022","start":"2022-06-16 15:00","match":"West Indies vs Bangladesh"},{"stadium":"R.Premadasa Stadium, Khettarama, Colombo","country":"Sri Lanka","region":"","tournament":"Sri Lanka vs Australia ODI Series 2022","start":"2022-06-21 10:00","match":"Sri Lanka vs Australia"}],"golf":[{"stadium":"Muirfield Village Gc","country":"United States of America","region":"","tournament":"the Memorial Tournament presented by Workday Round 3","start":"2022-06-04 12:35","match":"Lucas Herbert, Chan Kim"},{"stadium":"Muirfield Village Gc","country":"United States of America","region":"","tournament":"the Memorial Tournament presented by Workday Round 3","start":"2022-06-04 12:45","match":"Ryan Moore, Adam Scott"},{"stadium":"Muirfield Village Gc","country":"United States of America","region":"","tournament":"the Memorial Tournament presented by Workday Round 3","start":"2022-06-04 12:55","match":"Aaron Rai, Adam Schenk"},{"stadium":"Muirfield Village Gc","country":"United States of America","region":"","tournament":"the Memorial Tournament presented by Workday Round 3","start":"2022-06-04 13:05","match":"Kramer Hickok, David Lingmerth"}]}'
$footballJsonSofia = '{"football":[{"stadium":"Fc Rapid Bucuresti","country":"Romania","region":"","tournament":"UEFA Nations League","start":"2022-06-14 19:45","match":"Romania vs Montenegro"},{"stadium":"Inter Turku , Fc Tps","country":"Finland","region":"","tournament":"Finnish Veikkausliga","start":"2022-06-22 16:00","match":"FC Inter Turku vs FC Lahti"},{"stadium":"Drogheda United","country":"Ireland","region":"","tournament":"League of Ireland Premier Division","start":"2022-06-24 19:45","match":"Drogheda vs Sligo Rovers"},{"stadium":"Shelbourne Fc","country":"Ireland","region":"","tournament":"League of Ireland Premier Division","start":"2022-06-24 19:45","match":"Shelbourne vs Dundalk"},{"stadium":"Ucd","country":"Ireland","region":"","tournament":"League of Ireland Premier Division","start":"2022-06-24 19:45","match":"U.C.D vs Derry City"},{"stadium":"Cork City Fc","country":"Ireland","region":"","tournament":"League of Ireland First Division","start":"2022-06-24 19:45","match":"Cork City vs Cobh Ramblers"},{"stadium":"Waterford United Fc","country":"Ireland","region":"","tournament":"League of Ireland First Division","start":"2022-06-24 19:45","match":"Waterford United vs Bray Wanderers"},{"stadium":"Finn Harps","country":"Ireland","region":"","tournament":"League of Ireland Premier Division","start":"2022-06-24 20:00","match":"Finn Harps vs St Patricks Athletic"},{"stadium":"Shamrock Rovers(A)","country":"Ireland","region":"","tournament":"League of Ireland Premier Division","start":"2022-06-24 20:00","match":"Shamrock Rovers vs Bohemians"},{"stadium":"Cincinnati Kings","country":"Sofia","region":"","tournament":"American MLS League","start":"2022-06-25 00:30","match":"FC Cincinnati vs Orlando City SC"}],"cricket":[{"stadium":"R.Premadasa Stadium, Khettarama, Colombo","country":"Sri Lanka","region":"","tournament":"Sri Lanka vs Australia ODI Series 2022","start":"2022-06-21 10:00","match":"Sri Lanka vs Australia"}],"golf":[{"stadium":"Muirfield Village Gc","country":"United States of America","region":"","tournament":"the Memorial Tournament presented by Workday Round 3","start":"2022-06-04 12:35","match":"Lucas Herbert, Chan Kim"},{"stadium":"Muirfield Village Gc","country":"United States of America","region":"","tournament":"the Memorial Tournament presented by Workday Round 3","start":"2022-06-04 12:45","match":"Ryan Moore, Adam Scott"},{"stadium":"Muirfield Village Gc","country":"United States of America","region":"","tournament":"the Memorial Tournament presented by Workday Round 3","start":"2022-06-04 12:55","match":"Aaron Rai, Adam Schenk"},{"stadium":"Muirfield Village Gc","country":"United States of America","region":"","tournament":"the Memorial Tournament presented by Workday Round 3","start":"2022-06-04 13:05","match":"Kramer Hickok, David Lingmerth"}]}'

$wsFootball.Range("A2").Value = 0
$wsFootball.Range("B2").Value = $listMon
$wsFootball.Range("C2").Value = $footballJsonBulgaria
$wsFootball.Range("D2").Value = $footballLabel

$wsFootball.Range("A3").Value = 1
$wsFootball.Range("B3").Value = $listTue
$wsFootball.Range("C3").Value = $footballJsonSofia
$wsFootball.Range("D3").Value = $footballLabel

$wsFootball.Range("A4").Value = 1
$wsFootball.Range("B4").Value = $listTue
$wsFootball.Range("C4").Value = $footballJsonSofia
$wsFootball.Range("D4").Value = $footballLabel

$wsFootball.Range("A5").Value = 1
$wsFootball.Range("B5").Value = $listTue
$wsFootball.Range("C5").Value = $footballJsonSofia
$wsFootball.Range("D5").Value = $footballLabel

$wsFootball.Range("H10").Select()

# ---------------------------------------------------------------------
# Selections / active sheet bookkeeping to mirror the saved workbook
# state (forecast ends up being the active tab).
# ---------------------------------------------------------------------
$wsCurrent.Range("F5").Select()
$wsForecast.Range("D1").Select()
$wsForecast.Activate()
